# vendors_list.xlsx — add newly-tracked vendors turned up while tuning the
# vendor-matching regex (see commit message), and touch up the sheet the
# way the author evidently did by hand afterwards (widen the name column,
# set a print paper size/orientation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the four new vendor names found (rows 8-11). Setting .Value on
# these cells grows sheetData/dimension and appends the matching entries
# to sharedStrings.xml automatically.
$ws.Range("A8").Value  = "MagLink"
$ws.Range("A9").Value  = "MikoPBX"
$ws.Range("A10").Value = "eosphoros-ai"
$ws.Range("A11").Value = "HKUDS"

# Leave the cursor on the last entered row, matching the saved selection.
$ws.Range("A11").Select() | Out-Null

# The vendor names are now longer, so the name column got manually widened.
$ws.Columns.Item(1).ColumnWidth = 10

# Page setup was touched as well (paper size / orientation recorded).
$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
